# Atualizado por script em 12-11-2023 14:45
#
# This workbook is a scraped betting-odds table (one row per match). The
# scraping run that produced this commit re-fetched the fixture list and:
#   1) a handful of matches that share the exact same kickoff date/time
#      (columns A-E are identical within each such block) came back from the
#      source in a different order, so their F:V ("home" .. "url_partida")
#      payloads need to be rotated/swapped into the new row order while A:E
#      (Indice/pais/torneio/temporada/data_partida) stay untouched;
#   2) two brand-new matches (rows 126 and 127) were appended at the end.
#
# Below: first snapshot the F:V payload of every row that participates in a
# reshuffle (so later writes never clobber a value still needed), then write
# each row's new payload, then append the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot current F:V ("home" .. "url_partida") payloads ---
# for every row involved in a same-date reshuffle.
$row3   = $ws.Range("F3:V3").Value2
$row4   = $ws.Range("F4:V4").Value2
$row5   = $ws.Range("F5:V5").Value2
$row6   = $ws.Range("F6:V6").Value2
$row7   = $ws.Range("F7:V7").Value2
$row8   = $ws.Range("F8:V8").Value2
$row55  = $ws.Range("F55:V55").Value2
$row56  = $ws.Range("F56:V56").Value2
$row74  = $ws.Range("F74:V74").Value2
$row76  = $ws.Range("F76:V76").Value2
$row77  = $ws.Range("F77:V77").Value2
$row78  = $ws.Range("F78:V78").Value2
$row81  = $ws.Range("F81:V81").Value2
$row82  = $ws.Range("F82:V82").Value2
$row83  = $ws.Range("F83:V83").Value2
$row84  = $ws.Range("F84:V84").Value2
$row85  = $ws.Range("F85:V85").Value2
$row86  = $ws.Range("F86:V86").Value2
$row91  = $ws.Range("F91:V91").Value2
$row92  = $ws.Range("F92:V92").Value2
$row93  = $ws.Range("F93:V93").Value2
$row94  = $ws.Range("F94:V94").Value2
$row105 = $ws.Range("F105:V105").Value2
$row106 = $ws.Range("F106:V106").Value2
$row107 = $ws.Range("F107:V107").Value2

# --- Step 2: write each row's new F:V payload (new row <- old row) ---
$ws.Range("F3:V3").Value2   = $row4
$ws.Range("F4:V4").Value2   = $row5
$ws.Range("F5:V5").Value2   = $row6
$ws.Range("F6:V6").Value2   = $row7
$ws.Range("F7:V7").Value2   = $row8
$ws.Range("F8:V8").Value2   = $row3

$ws.Range("F55:V55").Value2 = $row56
$ws.Range("F56:V56").Value2 = $row55

$ws.Range("F74:V74").Value2 = $row76
$ws.Range("F76:V76").Value2 = $row74

$ws.Range("F77:V77").Value2 = $row78
$ws.Range("F78:V78").Value2 = $row77

$ws.Range("F81:V81").Value2 = $row86
$ws.Range("F82:V82").Value2 = $row85
$ws.Range("F83:V83").Value2 = $row84
$ws.Range("F84:V84").Value2 = $row82
$ws.Range("F85:V85").Value2 = $row81
$ws.Range("F86:V86").Value2 = $row83

$ws.Range("F91:V91").Value2 = $row93
$ws.Range("F92:V92").Value2 = $row94
$ws.Range("F93:V93").Value2 = $row91
$ws.Range("F94:V94").Value2 = $row92

$ws.Range("F105:V105").Value2 = $row106
$ws.Range("F106:V106").Value2 = $row107
$ws.Range("F107:V107").Value2 = $row105

# --- Step 3: append the two new matches as rows 126 and 127 ---
# Column A (Indice) is bold + bordered (style of existing A125) and column E
# (data_partida) carries the date/time number format (style of existing
# E125); copy those formats over via PasteSpecial before writing the values
# so the new cells reuse the workbook's existing style indices instead of
# minting new ones.

# Row 126: Samorin 2-3 Povazska Bystrica
$ws.Range("A125").Copy()
$ws.Range("A126").PasteSpecial(-4122)
$ws.Range("A126").Value2 = 125
$ws.Range("B126").Value2 = "slovakia"
$ws.Range("C126").Value2 = "2-liga"
$ws.Range("D126").Value2 = "2023-2024"
$ws.Range("E125").Copy()
$ws.Range("E126").PasteSpecial(-4122)
$ws.Range("E126").Value2 = 45242.4375

$arr126 = New-Object 'object[,]' 1,17
$arr126[0,0]  = "Samorin"
$arr126[0,1]  = 2
$arr126[0,2]  = "Povazska Bystrica"
$arr126[0,3]  = 3
$arr126[0,4]  = 2.7
$arr126[0,5]  = "10/11/2023 04:42"
$arr126[0,6]  = 2.65
$arr126[0,7]  = "12/11/2023 10:19"
$arr126[0,8]  = 3.34
$arr126[0,9]  = "10/11/2023 04:42"
$arr126[0,10] = 3.75
$arr126[0,11] = "12/11/2023 10:19"
$arr126[0,12] = 2.3
$arr126[0,13] = "10/11/2023 04:42"
$arr126[0,14] = 2.34
$arr126[0,15] = "12/11/2023 10:19"
$arr126[0,16] = "https://www.betexplorer.com/football/slovakia/2-liga/samorin-povazska-bystrica/IqaSYXlk/"
$ws.Range("F126:V126").Value2 = $arr126

# Row 127: Zilina B 3-0 D. Kubin
$ws.Range("A125").Copy()
$ws.Range("A127").PasteSpecial(-4122)
$ws.Range("A127").Value2 = 126
$ws.Range("B127").Value2 = "slovakia"
$ws.Range("C127").Value2 = "2-liga"
$ws.Range("D127").Value2 = "2023-2024"
$ws.Range("E125").Copy()
$ws.Range("E127").PasteSpecial(-4122)
$ws.Range("E127").Value2 = 45242.4375

$arr127 = New-Object 'object[,]' 1,17
$arr127[0,0]  = "Zilina B"
$arr127[0,1]  = 3
$arr127[0,2]  = "D. Kubin"
$arr127[0,3]  = 0
$arr127[0,4]  = 1.48
$arr127[0,5]  = "10/11/2023 04:42"
$arr127[0,6]  = 1.24
$arr127[0,7]  = "12/11/2023 10:20"
$arr127[0,8]  = 4.41
$arr127[0,9]  = "10/11/2023 04:42"
$arr127[0,10] = 6.41
$arr127[0,11] = "12/11/2023 10:20"
$arr127[0,12] = 4.65
$arr127[0,13] = "10/11/2023 04:42"
$arr127[0,14] = 8.76
$arr127[0,15] = "12/11/2023 10:20"
$arr127[0,16] = "https://www.betexplorer.com/football/slovakia/2-liga/zilina-d-kubin/n18dvfRR/"
$ws.Range("F127:V127").Value2 = $arr127
